# Applies the "implemented CR, CP, F, AR metrics and generated consolidated report" update:
#  - Summary sheet: adds 9 new test-run rows (runs 21-29) above the TOTALS row,
#    and refreshes the TOTALS row with the new aggregate figures.
#  - Question Failure Rates sheet: refreshes total_runs/passed/failed/na/failure_rate
#    for every question now that there are 9 runs instead of 8.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

# Insert 9 blank rows right above the old TOTALS row (row 10) so the TOTALS
# row ends up at row 19, and rows 10-18 are free for the new run data.
$wsSummary.Range("A10:A18").EntireRow.Insert()

# New run rows: run, total_tests, passed, failed, na_not_evaluated, pass_rate,
# avg_precision_score, failed_questions, na_questions, execution_time_seconds
$newRuns = @(
    @(10, 21, 35, 25, 4, 6, "71.4%", 0.823,  "L, M, R, Q28",       "C, O, P, S, Q29, Q31"),
    @(11, 22, 35, 25, 4, 6, "71.4%", 0.8621, "L, M, R, Q27",       "B, C, O, P, S, Q31"),
    @(12, 23, 35, 27, 2, 6, "77.1%", 0.8996, "M, R",               "C, O, P, S, Q28, Q31"),
    @(13, 24, 35, 26, 3, 6, "74.3%", 0.8768, "L, M, R",            "C, O, P, S, Q28, Q31"),
    @(14, 25, 35, 26, 2, 7, "74.3%", 0.8988, "B, R",               "A, C, O, P, S, Q28, Q31"),
    @(15, 26, 35, 23, 4, 8, "65.7%", 0.8519, "D, M, R, Q27",       "A, C, O, P, S, Q28, Q29, Q31"),
    @(16, 27, 35, 24, 3, 8, "68.6%", 0.8551, "N, R, Q27",          "A, C, L, P, S, Q28, Q29, Q31"),
    @(17, 28, 35, 21, 6, 8, "60.0%", 0.8243, "D, F, L, M, N, R",   "A, C, O, P, S, Q28, Q29, Q31"),
    @(18, 29, 35, 21, 4, 10, "60.0%", 0.8733, "D, L, R, Q27",      "A, C, O, P, S, Q28, Q29, Q31, Q33, Q35")
)

$execTimes = @{
    10 = 199.41; 11 = 205.05; 12 = 198.58; 13 = 179.62; 14 = 176.8;
    15 = 167.67; 16 = 200.44; 17 = 182.98; 18 = 148.43
}

foreach ($row in $newRuns) {
    $r = $row[0]
    $wsSummary.Cells.Item($r, 1).Value = $row[1]
    $wsSummary.Cells.Item($r, 2).Value = $row[2]
    $wsSummary.Cells.Item($r, 3).Value = $row[3]
    $wsSummary.Cells.Item($r, 4).Value = $row[4]
    $wsSummary.Cells.Item($r, 5).Value = $row[5]
    Set-TextValue $wsSummary.Cells.Item($r, 6) $row[6]
    $wsSummary.Cells.Item($r, 7).Value = $row[7]
    $wsSummary.Cells.Item($r, 8).Value = $row[8]
    $wsSummary.Cells.Item($r, 9).Value = $row[9]
    $wsSummary.Cells.Item($r, 10).Value = $execTimes[$r]
}

# Refresh the TOTALS row (now row 19) with the consolidated totals.
$wsSummary.Cells.Item(19, 1).Value = "TOTALS"
$wsSummary.Cells.Item(19, 2).Value = "17 runs"
$wsSummary.Cells.Item(19, 3).Value = ""
$wsSummary.Cells.Item(19, 4).Value = ""
$wsSummary.Cells.Item(19, 5).Value = 119
Set-TextValue $wsSummary.Cells.Item(19, 6) "69.2%"
$wsSummary.Cells.Item(19, 7).Value = 0.86
$wsSummary.Cells.Item(19, 8).Value = ""
$wsSummary.Cells.Item(19, 9).Value = ""
$wsSummary.Cells.Item(19, 10).Value = 861.01

# ---------------------------------------------------------------------------
# 2. Question Failure Rates sheet
# ---------------------------------------------------------------------------
$wsRates = $wb.Worksheets.Item("Question Failure Rates")

# row -> total_runs, passed, failed, na, failure_rate
$rateUpdates = @{
    2  = @(9, 4, 0, 5, "0.0%")
    3  = @(9, 7, 1, 1, "11.1%")
    4  = @(9, 0, 0, 9, "0.0%")
    5  = @(9, 6, 3, 0, "33.3%")
    6  = @(9, 9, 0, 0, "0.0%")
    7  = @(9, 8, 1, 0, "11.1%")
    8  = @(9, 9, 0, 0, "0.0%")
    9  = @(9, 9, 0, 0, "0.0%")
    10 = @(9, 9, 0, 0, "0.0%")
    11 = @(9, 9, 0, 0, "0.0%")
    12 = @(9, 9, 0, 0, "0.0%")
    13 = @(9, 3, 5, 1, "55.6%")
    14 = @(9, 3, 6, 0, "66.7%")
    15 = @(9, 7, 2, 0, "22.2%")
    16 = @(9, 1, 0, 8, "0.0%")
    17 = @(9, 0, 0, 9, "0.0%")
    18 = @(9, 9, 0, 0, "0.0%")
    19 = @(9, 0, 9, 0, "100.0%")
    20 = @(9, 0, 0, 9, "0.0%")
    21 = @(9, 9, 0, 0, "0.0%")
    22 = @(9, 9, 0, 0, "0.0%")
    23 = @(9, 9, 0, 0, "0.0%")
    24 = @(9, 9, 0, 0, "0.0%")
    25 = @(9, 9, 0, 0, "0.0%")
    26 = @(9, 9, 0, 0, "0.0%")
    27 = @(9, 9, 0, 0, "0.0%")
    28 = @(9, 5, 4, 0, "44.4%")
    29 = @(9, 1, 1, 7, "11.1%")
    30 = @(9, 4, 0, 5, "0.0%")
    31 = @(9, 9, 0, 0, "0.0%")
    32 = @(9, 0, 0, 9, "0.0%")
    33 = @(9, 9, 0, 0, "0.0%")
    34 = @(9, 8, 0, 1, "0.0%")
    35 = @(9, 9, 0, 0, "0.0%")
    36 = @(9, 8, 0, 1, "0.0%")
}

foreach ($r in $rateUpdates.Keys) {
    $vals = $rateUpdates[$r]
    $wsRates.Cells.Item($r, 3).Value = $vals[0]
    $wsRates.Cells.Item($r, 4).Value = $vals[1]
    $wsRates.Cells.Item($r, 5).Value = $vals[2]
    $wsRates.Cells.Item($r, 6).Value = $vals[3]
    Set-TextValue $wsRates.Cells.Item($r, 7) $vals[4]
}

Write-Host "Summary and Question Failure Rates sheets updated."
